$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.428.21"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.067.50"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.73"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.13"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0775"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.371.05"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.35"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.66"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.778"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.066.38"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.323.44"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.21"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.51"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.35"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.15"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  +4.54%  "
$ws.Range("E29").Value = "  -6.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.08"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.55"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  -4.31%  "
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.492.50"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0952"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.06"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.20"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.24"
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.257.61"
$ws.Range("E51").Value = "  -0.01%  "
